$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells stay text (preserve leading zeros, keep as string type) like the
# original inlineStr cells.
$ws.Range("A1:G3").NumberFormat = "@"

# Update rows 1-3 with new subject/teacher/time data
$ws.Range("A1").Value = "49971"
$ws.Range("B1").Value = "Vision y Animacion por Comput."
$ws.Range("C1").Value = "JUAREZ - PEREZ SILVESTRE"
$ws.Range("D1").Value = "1000"
$ws.Range("E1").Value = "1059"

$ws.Range("A2").Value = "49971"
$ws.Range("B2").Value = "Vision y Animacion por Comput."
$ws.Range("C2").Value = "JUAREZ - PEREZ SILVESTRE"
$ws.Range("D2").Value = "0900"
$ws.Range("E2").Value = "1059"

$ws.Range("A3").Value = "49971"
$ws.Range("B3").Value = "Vision y Animacion por Comput."
$ws.Range("C3").Value = "JUAREZ - PEREZ SILVESTRE"
$ws.Range("D3").Value = "0900"
$ws.Range("E3").Value = "1059"
$ws.Range("G3").Value = "1CCO4/308"

# Remove rows 4 through 9 (old Sistemas Operativos I / Web Semantica entries)
$ws.Range("A4:G9").Delete()
